# E5 initial test samples and blue tank titrations 0504 dmb
#
# The CRM (certified reference material) bottle used since 20210418 was
# re-weighed, correcting the prior C31 "CRM value" entry, and a new titration
# row (20210504) was appended below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: correct the CRM value / batch # for the 20210428 sample ---
$ws.Range("C31").Value = 2224.4699999999998
$ws.Range("E31").Value = 180

# --- Row 32: new titration sample, 20210504 ---
$ws.Range("A32").Value = 20210504
$ws.Range("B32").Value = 2224.8229999999999
$ws.Range("C32").Value = 2224.4699999999998
$ws.Range("D32").Formula = "=100*(B32-C32)/C32"
$ws.Range("E32").Value = 180
$ws.Range("F32").Value = "CRM opened 20210418"

# Leave the selection where the author ended up after entering the new row
$ws.Range("I29").Select()
